$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-10: records rotated (row10 -> row3, row3->row4, row4->row5, ...,
#     row9->row10). Only the record-specific columns change; shared
#     per-sighting columns (C, I, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT,
#     AW, AX, AY) stay identical across all these rows and are unaffected.

# Row 3 (was row 10's data)
$ws.Range("A3").Value = 67898328
$ws.Range("B3").Value = 73693
$ws.Range("E3").Value = 6440
$ws.Range("F3").Value = "Vitgrynig nållav"
$ws.Range("G3").Value = "Chaenotheca subroscida"
$ws.Range("H3").Value = "(Eitner) Zahlbr."
$ws.Range("P3").Value = "Naturskog SV Ladumyrberget, syd kraftledningen, Ång"
$ws.Range("Q3").Value = 584489.8304919918
$ws.Range("R3").Value = 7048424.423746439

# Row 4 (was row 3's data)
$ws.Range("A4").Value = 67892102
$ws.Range("B4").Value = 77541
$ws.Range("E4").Value = 185
$ws.Range("F4").Value = "Violettgrå tagellav"
$ws.Range("G4").Value = "Bryoria nadvornikiana"
$ws.Range("H4").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q4").Value = 584742.9303808049
$ws.Range("R4").Value = 7048428.485565316

# Row 5 (was row 4's data)
$ws.Range("A5").Value = 67892121
$ws.Range("B5").Value = 77259
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 228912
$ws.Range("F5").Value = "Mörk kolflarnlav"
$ws.Range("G5").Value = "Carbonicola myrmecina"
$ws.Range("H5").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q5").Value = 584700.4782698094
$ws.Range("R5").Value = 7048463.032207629
$ws.Range("S5").Value = 5

# Row 6 (was row 5's data)
$ws.Range("A6").Value = 67892107
$ws.Range("B6").Value = 76863
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 498
$ws.Range("F6").Value = "Liten sotlav"
$ws.Range("G6").Value = "Acolium karelicum"
$ws.Range("H6").Value = "(Vain.) M.Prieto & Wedin"
$ws.Range("S6").Value = 25

# Row 7 (was row 6's data)
$ws.Range("A7").Value = 67892103
$ws.Range("B7").Value = 77506
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 584742.9303808049
$ws.Range("R7").Value = 7048428.485565316

# Row 8 (was row 7's data)
$ws.Range("A8").Value = 67892123
$ws.Range("B8").Value = 77258
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6446
$ws.Range("F8").Value = "Kolflarnlav"
$ws.Range("G8").Value = "Carbonicola anthracophila"
$ws.Range("H8").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q8").Value = 584700.4782698094
$ws.Range("R8").Value = 7048463.032207629
$ws.Range("S8").Value = 5

# Row 9 (was row 8's data)
$ws.Range("A9").Value = 67892109
$ws.Range("B9").Value = 78596
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 6462
$ws.Range("F9").Value = "Stuplav"
$ws.Range("G9").Value = "Nephroma bellum"
$ws.Range("H9").Value = "(Spreng.) Tuck."

# Row 10 (was row 9's data)
$ws.Range("A10").Value = 67892108
$ws.Range("B10").Value = 78569
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = "Lunglav"
$ws.Range("G10").Value = "Lobaria pulmonaria"
$ws.Range("H10").Value = "(L.) Hoffm."
$ws.Range("P10").Value = "Naturskog SV Ladumyrberget, Ång"
$ws.Range("Q10").Value = 584742.9303808049
$ws.Range("R10").Value = 7048428.485565316
$ws.Range("S10").Value = 25

# --- Rows 26-27: the two sighting records swap places (with the
#     Taxonsorteringsordning in column B receiving fresh values rather than
#     simply swapping).

# Row 26 (now holds the "Knärot / Goodyera repens" record)
$ws.Range("A26").Value = 111935024
$ws.Range("B26").Value = 96735
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = "Knärot"
$ws.Range("G26").Value = "Goodyera repens"
$ws.Range("H26").Value = "(L.) R. Br."
$ws.Range("I26").Value = "10"
$ws.Range("P26").Value = "Sollefteå (Sollefteå), Ång"
$ws.Range("Q26").Value = 584598
$ws.Range("R26").Value = 7048260
$ws.Range("S26").Value = 2
$ws.Range("Z26").Value = "19:51"
$ws.Range("AB26").Value = "19:51"
$ws.Range("AW26").Value = "Kim Hultgren"
$ws.Range("AX26").Value = "Kim Hultgren"

# Row 27 (now holds the "Kolflarnlav / Carbonicola anthracophila" record)
$ws.Range("A27").Value = 111934989
$ws.Range("B27").Value = 77402
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 6446
$ws.Range("F27").Value = "Kolflarnlav"
$ws.Range("G27").Value = "Carbonicola anthracophila"
$ws.Range("H27").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("I27").Value = ""
$ws.Range("P27").Value = "Skumsåstjärnen (Skumsåstjärnen), Ång"
$ws.Range("Q27").Value = 584600
$ws.Range("R27").Value = 7048296
$ws.Range("S27").Value = 1
$ws.Range("Z27").Value = ""
$ws.Range("AB27").Value = ""
$ws.Range("AW27").Value = "Kamilla Andersson"
$ws.Range("AX27").Value = "Kamilla Andersson"
